# Reorder the title/timestamp/uri columns for rows 2-13 of the
# historical_distance sheet (time-bucket analysis refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Presidential Ratings"
$ws.Cells.Item(2, 2).Value = "1-01-01T00:00:00UTC"
$ws.Cells.Item(2, 5).Value = "https://insideelections.com/ratings/president"

$ws.Cells.Item(3, 1).Value = "Iowa: Policy Priorities and the Election"
$ws.Cells.Item(3, 2).Value = "2020-10-22T00:00:00UTC"
$ws.Cells.Item(3, 5).Value = "https://www.americanactionforum.org/insight/iowa-policy-priorities-and-the-election-october-update/"

$ws.Cells.Item(4, 1).Value = "2020 President - Sabato's Crystal Ball"
$ws.Cells.Item(4, 2).Value = "1-01-01T00:00:00UTC"
$ws.Cells.Item(4, 5).Value = "http://centerforpolitics.org/crystalball/2020-president/"

$ws.Cells.Item(5, 1).Value = "Libertarian's announce caucus results"
$ws.Cells.Item(5, 2).Value = "2020-02-10T15:23:00UTC"
$ws.Cells.Item(5, 5).Value = "https://www.journalexpress.net/news/local_news/libertarian-s-announce-caucus-results/article_95ef52f4-4c4b-11ea-9490-8f4b504e69d2.html"

$ws.Cells.Item(6, 1).Value = "Low voter turnout at the Iowa Libertarian Party Caucus"
$ws.Cells.Item(6, 2).Value = "2020-02-10T14:09:00UTC"
$ws.Cells.Item(6, 5).Value = "https://www.oskaloosa.com/iowa/low-voter-turnout-at-the-iowa-libertarian-party-caucus/article_4b2fa8cc-4c41-11ea-91ab-97b3e4d94934.html"

$ws.Cells.Item(7, 1).Value = "Buttigieg, Sanders in near tie"
$ws.Cells.Item(7, 2).Value = "2020-02-04T00:00:00UTC"
$ws.Cells.Item(7, 5).Value = "https://www.usatoday.com/story/news/politics/elections/2020/02/06/iowa-caucus-results-pete-buttigieg-bernie-sanders-close-delegate-race/4675289002/"

$ws.Cells.Item(8, 1).Value = "Iowa Democratic Party Announces Delegation to National Convention"
$ws.Cells.Item(8, 2).Value = "2020-06-13T13:15:51UTC"
$ws.Cells.Item(8, 5).Value = "https://iowademocrats.org/iowa-democratic-party-announces-delegation-national-convention/"

$ws.Cells.Item(9, 1).Value = "Whoever Wins Iowa, They Won’t Be Back"
$ws.Cells.Item(9, 2).Value = "2020-02-03T05:05:00UTC"
$ws.Cells.Item(9, 5).Value = "https://www.politico.com/news/magazine/2020/02/03/whoever-wins-iowa-wont-be-back-110439"

$ws.Cells.Item(10, 1).Value = "Voting & Elections Toolkits"
$ws.Cells.Item(10, 2).Value = "1-01-01T00:00:00UTC"
$ws.Cells.Item(10, 5).Value = "https://godort.libguides.com/votingtoolkit/texas"

$ws.Cells.Item(11, 1).Value = "Biden dominates the electoral map, but here's how the race could tighten"
$ws.Cells.Item(11, 2).Value = "2020-08-06T13:13:00UTC"
$ws.Cells.Item(11, 5).Value = "https://www.nbcnews.com/politics/meet-the-press/biden-dominates-electoral-map-here-s-how-race-could-tighten-n1236001"

$ws.Cells.Item(12, 1).Value = "Heartland Poll Release: Biden Leads in Midwest"
$ws.Cells.Item(12, 2).Value = "2020-08-06T15:46:54UTC"
$ws.Cells.Item(12, 5).Value = "https://www.focusonruralamerica.com/2020/08/06/heartland-poll-biden-leads-in-midwest/"

$ws.Cells.Item(13, 1).Value = "2020 Electoral Interactive Map"
$ws.Cells.Item(13, 2).Value = "1-01-01T00:00:00UTC"
$ws.Cells.Item(13, 5).Value = "https://abcnews.go.com/Politics/2020-Electoral-Interactive-Map?basemap=71662160&promoref=brandpromo"
